# Inserts a new data row for "Macroferia Regional de Talca - Brócoli" at row 367,
# shifting all existing rows 367:493 down to 368:494 (dimension grows from
# A1:R493 to A1:R494), then populates the newly-inserted row 367 with the
# new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 367, pushing everything below
# (including the former row 367) down by one.
$ws.Rows(367).Insert()

# Populate the new row 367 with the new record.
$ws.Range("A367").Value = 5
$ws.Range("B367").Value = "Macroferia Regional de Talca"
$ws.Range("C367").Value = "Maule"
$ws.Range("D367").Value = 44985
$ws.Range("E367").Value = 7
$ws.Range("F367").Value = 100112023
$ws.Range("G367").Value = "Brócoli"
$ws.Range("H367").Value = "Sin especificar"
$ws.Range("I367").Value = "Primera"
$ws.Range("J367").Value = 4000
$ws.Range("K367").Value = 700
$ws.Range("L367").Value = 700
$ws.Range("M367").Value = 700
$ws.Range("N367").Value = "$/unidad"
$ws.Range("O367").Value = "Región del Maule"
$ws.Range("P367").Value = 700
$ws.Range("Q367").Value = 1
$ws.Range("R367").Value = "Hortaliza"
